# International Conflict deck: append entity-type numeric IDs to labels,
# fill in the "Overview" slide's placeholder title/subtitle, and shrink the
# (now much narrower) Overview body placeholder box.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Overview") - body placeholder (shape 2):
#   - narrow the box from 8520600 EMU wide to 1678200 EMU wide
#   - the single (until now empty) run becomes "International Conflict (10277)"
# ---------------------------------------------------------------------------
$sOverview = $p.Slides.Item(2)
$shOverviewBody = $sOverview.Shapes.Item(2)

$shOverviewBody.Width = 1678200 / 12700

$trOverviewBody = $shOverviewBody.TextFrame.TextRange
$trOverviewBody.Text = "International Conflict (10277)"
$trOverviewBody.LanguageID = "en"

# ---------------------------------------------------------------------------
# Slide 3 ("Battle Episode")
# ---------------------------------------------------------------------------
$sBattle = $p.Slides.Item(3)

$shBattleTitle = $sBattle.Shapes.Item(1)
$shBattleTitle.TextFrame.TextRange.Text = "Battle Episode (10286)"

$shBattleBody = $sBattle.Shapes.Item(2)
$trBattleBody = $shBattleBody.TextFrame.TextRange

# Each bullet paragraph starts with a bold "label" run; Runs() merges the
# trailing same-formatted runs of each paragraph into a single logical run,
# so the label is always the odd-numbered entry below.
$trBattleBody.Runs(1, 1).Text = "Attack (10278)"
$trBattleBody.Runs(3, 1).Text = "Scarcity (10282)"
$trBattleBody.Runs(5, 1).Text = "Injury (10280)"
$trBattleBody.Runs(7, 1).Text = "Death (10279)"
$trBattleBody.Runs(9, 1).Text = "Damage (10281)"
$trBattleBody.Runs(11, 1).Text = "Transport (10283)"
$trBattleBody.Runs(13, 1).Text = "Demonstration (10284)"
$trBattleBody.Runs(15, 1).Text = "Arrest (10285)"

# ---------------------------------------------------------------------------
# Slide 4 ("Ceasefire Episode")
# ---------------------------------------------------------------------------
$sCeasefire = $p.Slides.Item(4)

$shCeasefireTitle = $sCeasefire.Shapes.Item(1)
$shCeasefireTitle.TextFrame.TextRange.Text = "Ceasefire Episode (10292)"

$shCeasefireBody = $sCeasefire.Shapes.Item(2)
$trCeasefireBody = $shCeasefireBody.TextFrame.TextRange

$trCeasefireBody.Runs(1, 1).Text = "Remote Communication (10287)"
$trCeasefireBody.Runs(3, 1).Text = "Negotiation (10288)"
$trCeasefireBody.Runs(5, 1).Text = "Negotiation Result (10289)"
$trCeasefireBody.Runs(7, 1).Text = "Rejection (10290)"
$trCeasefireBody.Runs(9, 1).Text = "Agreement (10291)"
$trCeasefireBody.Runs(11, 1).Text = "Announcement (10293)"
$trCeasefireBody.Runs(13, 1).Text = "Withdrawal (10294)"
